$d = $word.ActiveDocument

function Update-OfficeHoursRoom($searchText, $newRoom) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host ("NOT FOUND: " + $searchText)
        return $null
    }
    $matchStart = $rng.Start
    $matchEnd = $rng.End
    $tbdLen = 3  # length of "TBD"

    # Step 1: split the "...Campbell " prefix away from the preceding run
    # (e.g. "M 4:00p -") by toggling Bold off/on - this forces a run
    # boundary without altering the resulting formatting.
    $prefixRange = $d.Range($matchStart, $matchEnd - $tbdLen)
    $prefixRange.Bold = 1
    $prefixRange.Bold = 0

    # Step 2: replace "TBD" with the new room number.
    $tbdRange = $d.Range($matchEnd - $tbdLen, $matchEnd)
    $tbdRange.Text = $newRoom

    # Step 3: split the new room-number text into its own run by toggling
    # Bold off/on again.
    $newLen = $newRoom.Length
    $newRange = $d.Range($matchEnd - $tbdLen, $matchEnd - $tbdLen + $newLen)
    $newRange.Bold = 1
    $newRange.Bold = 0

    return $newRange
}

Update-OfficeHoursRoom " 4:30p Campbell TBD" "121" | Out-Null
Update-OfficeHoursRoom " 12:30p Campbell TBD" "233" | Out-Null
$lastRange = Update-OfficeHoursRoom " 2:30p Campbell TBD" "233"

# Move the _GoBack bookmark from the empty paragraph further down to right
# after the newly-typed "233".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$endPoint = $d.Range($lastRange.End, $lastRange.End)
$d.Bookmarks.Add("_GoBack", $endPoint) | Out-Null
